$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 23.34546866666667
$ws.Range("H2").Value = 70.036406
$ws.Range("I2").Value = 0.4715073400272545
$ws.Range("J2").Value = 0.4715073400272545
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 223.2367336666667
$ws.Range("N2").Value = 669.710201
$ws.Range("O2").Value = 0.9523995969492647
$ws.Range("P2").Value = 0.9523995969492646
$ws.Range("Q2").Value = 5211.566171064178
$ws.Range("R2").Value = 46904.0955395776
$ws.Range("S2").Value = 0.449063400600577
$ws.Range("T2").Value = 0.449063400600577

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 23.34546866666667
$ws.Range("H3").Value = 70.036406
$ws.Range("I3").Value = 0.4715073400272545
$ws.Range("J3").Value = 0.4715073400272545
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.280784333333333
$ws.Range("N3").Value = 12.842353
$ws.Range("O3").Value = 0.01826320071400582
$ws.Range("P3").Value = 0.01826320071400582
$ws.Range("Q3").Value = 99.93691652259088
$ws.Range("R3").Value = 899.4322487033179
$ws.Range("S3").Value = 0.008611233189044738
$ws.Range("T3").Value = 0.008611233189044738

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 23.34546866666667
$ws.Range("H4").Value = 70.036406
$ws.Range("I4").Value = 0.4715073400272545
$ws.Range("J4").Value = 0.4715073400272545
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.876463666666666
$ws.Range("N4").Value = 20.629391
$ws.Range("O4").Value = 0.02933720233672951
$ws.Range("P4").Value = 0.0293372023367295
$ws.Range("Q4").Value = 160.5342670676384
$ws.Range("R4").Value = 1444.808403608746
$ws.Range("S4").Value = 0.01383270623763268
$ws.Range("T4").Value = 0.01383270623763268

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.57337833333333
$ws.Range("H5").Value = 37.720135
$ws.Range("I5").Value = 0.2539439348061199
$ws.Range("J5").Value = 0.2539439348061199
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 223.2367336666667
$ws.Range("N5").Value = 669.710201
$ws.Range("O5").Value = 0.9523995969492647
$ws.Range("P5").Value = 0.9523995969492646
$ws.Range("Q5").Value = 2806.83991028857
$ws.Range("R5").Value = 25261.55919259713
$ws.Range("S5").Value = 0.2418561011570589
$ws.Range("T5").Value = 0.2418561011570589

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.57337833333333
$ws.Range("H6").Value = 37.720135
$ws.Range("I6").Value = 0.2539439348061199
$ws.Range("J6").Value = 0.2539439348061199
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.280784333333333
$ws.Range("N6").Value = 12.842353
$ws.Range("O6").Value = 0.01826320071400582
$ws.Range("P6").Value = 0.01826320071400582
$ws.Range("Q6").Value = 53.8239209864061
$ws.Range("R6").Value = 484.415288877655
$ws.Range("S6").Value = 0.004637829051468575
$ws.Range("T6").Value = 0.004637829051468575

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.57337833333333
$ws.Range("H7").Value = 37.720135
$ws.Range("I7").Value = 0.2539439348061199
$ws.Range("J7").Value = 0.2539439348061199
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.876463666666666
$ws.Range("N7").Value = 20.629391
$ws.Range("O7").Value = 0.02933720233672951
$ws.Range("P7").Value = 0.0293372023367295
$ws.Range("Q7").Value = 86.46037927642054
$ws.Range("R7").Value = 778.1434134877849
$ws.Range("S7").Value = 0.007450004597592385
$ws.Range("T7").Value = 0.007450004597592384

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 13.59357133333333
$ws.Range("H8").Value = 40.780714
$ws.Range("I8").Value = 0.2745487251666257
$ws.Range("J8").Value = 0.2745487251666257
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 223.2367336666667
$ws.Range("N8").Value = 669.710201
$ws.Range("O8").Value = 0.9523995969492647
$ws.Range("P8").Value = 0.9523995969492646
$ws.Range("Q8").Value = 3034.584463318168
$ws.Range("R8").Value = 27311.26016986351
$ws.Range("S8").Value = 0.2614800951916287
$ws.Range("T8").Value = 0.2614800951916287

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 13.59357133333333
$ws.Range("H9").Value = 40.780714
$ws.Range("I9").Value = 0.2745487251666257
$ws.Range("J9").Value = 0.2745487251666257
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.280784333333333
$ws.Range("N9").Value = 12.842353
$ws.Range("O9").Value = 0.01826320071400582
$ws.Range("P9").Value = 0.01826320071400582
$ws.Range("Q9").Value = 58.19114719778243
$ws.Range("R9").Value = 523.7203247800419
$ws.Range("S9").Value = 0.005014138473492505
$ws.Range("T9").Value = 0.005014138473492505

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.59357133333333
$ws.Range("H10").Value = 40.780714
$ws.Range("I10").Value = 0.2745487251666257
$ws.Range("J10").Value = 0.2745487251666257
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.876463666666666
$ws.Range("N10").Value = 20.629391
$ws.Range("O10").Value = 0.02933720233672951
$ws.Range("P10").Value = 0.0293372023367295
$ws.Range("Q10").Value = 93.4756993739082
$ws.Range("R10").Value = 841.2812943651738
$ws.Range("S10").Value = 0.008054491501504439
$ws.Range("T10").Value = 0.008054491501504437
